$wb = $excel.ActiveWorkbook

# --- 1. Update status text: "Ready for handoff" -> "In Translation" ---
# This string is shared across the Overview sheet (columns E/F, row 2)
# and the per-locale sheets (column C, row 2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Narrow the "Status" columns (Overview E:F, zh-cn C, de-de C) ---
# Target stored width is ~13.41 chars; the COM ColumnWidth setter here
# snaps to a 1/6-char grid, so 12.5 is the closest input that lands on
# the nearest achievable stored width.
$wsOverview.Range("E1:F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
